$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 41900.64
$ws.Range("J40").Value = 2074.625
$ws.Range("L40").Value = 2074.625
$ws.Range("N40").Value = -2424.625
$ws.Range("H64").Value = 86274.914
$ws.Range("I64").Value = 202400
$ws.Range("J64").Value = 3328.4285
$ws.Range("K64").Value = 202400
$ws.Range("L64").Value = 3328.4285
$ws.Range("M64").Value = -202152
$ws.Range("N64").Value = -3824.4285
$ws.Range("H67").Value = 86274.914
$ws.Range("I67").Value = 202400
$ws.Range("J67").Value = 3328.4285
$ws.Range("K67").Value = 202400
$ws.Range("L67").Value = 3328.4285
$ws.Range("M67").Value = -201542
$ws.Range("N67").Value = -5044.4285
$ws.Range("H76").Value = 4181.4
$ws.Range("I76").Value = 3002
$ws.Range("J76").Value = 4476.25
$ws.Range("K76").Value = 3002
$ws.Range("L76").Value = 4476.25
$ws.Range("M76").Value = -2687
$ws.Range("N76").Value = -5106.25
$ws.Range("H79").Value = 4181.4
$ws.Range("I79").Value = 3002
$ws.Range("J79").Value = 4476.25
$ws.Range("K79").Value = 3002
$ws.Range("L79").Value = 4476.25
$ws.Range("M79").Value = -1910
$ws.Range("N79").Value = -6660.25
$ws.Range("H80").Value = 111599.72
$ws.Range("I80").Value = 143002.42
$ws.Range("J80").Value = 91616.17999999999
$ws.Range("K80").Value = 429007.26
$ws.Range("L80").Value = 274848.54
$ws.Range("M80").Value = -428009.26
$ws.Range("N80").Value = -276844.54
$ws.Range("H83").Value = 111599.72
$ws.Range("I83").Value = 143002.42
$ws.Range("J83").Value = 91616.17999999999
$ws.Range("K83").Value = 1287021.78
$ws.Range("L83").Value = 824545.6199999999
$ws.Range("M83").Value = -1282029.78
$ws.Range("N83").Value = -834529.6199999999
$ws.Range("H137").Value = 1712.6111
$ws.Range("I137").Value = 1728.9166
$ws.Range("J137").Value = 1680
$ws.Range("K137").Value = 5186.7498
$ws.Range("L137").Value = 5040
$ws.Range("M137").Value = -2636.7498
$ws.Range("N137").Value = -10140
$ws.Range("H138").Value = 2337.407
$ws.Range("I138").Value = 2228.5334
$ws.Range("J138").Value = 2360.4084
$ws.Range("K138").Value = 6685.600199999999
$ws.Range("L138").Value = 7081.225199999999
$ws.Range("M138").Value = -1545.600199999999
$ws.Range("N138").Value = -17361.2252

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1180.9459
$ws.Range("I61").Value = 1128.4193
$ws.Range("J61").Value = 1452.3334
$ws.Range("K61").Value = 1128.4193
$ws.Range("L61").Value = 1452.3334
$ws.Range("M61").Value = -916.4193
$ws.Range("N61").Value = -1876.3334
$ws.Range("H74").Value = 2046.4688
$ws.Range("I74").Value = 1027.76
$ws.Range("J74").Value = 5684.7144
$ws.Range("K74").Value = 1027.76
$ws.Range("L74").Value = 5684.7144
$ws.Range("M74").Value = -153.76
$ws.Range("N74").Value = -7432.7144
$ws.Range("H77").Value = 2046.4688
$ws.Range("I77").Value = 1027.76
$ws.Range("J77").Value = 5684.7144
$ws.Range("K77").Value = 5138.8
$ws.Range("L77").Value = 28423.572
$ws.Range("M77").Value = -770.8000000000002
$ws.Range("N77").Value = -37159.572
$ws.Range("H88").Value = 3512
$ws.Range("I88").Value = 3442.2856
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 3442.2856
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -3036.2856
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 3512
$ws.Range("I91").Value = 3442.2856
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 3442.2856
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -2038.2856
$ws.Range("N91").Value = -6808
$ws.Range("H98").Value = 10564.2
$ws.Range("J98").Value = 10564.2
$ws.Range("L98").Value = 10564.2
$ws.Range("N98").Value = -16554.2
$ws.Range("H132").Value = 2673.1333
$ws.Range("I132").Value = 2610.3
$ws.Range("K132").Value = 7830.900000000001
$ws.Range("M132").Value = -5300.900000000001
$ws.Range("H136").Value = 1180.9459
$ws.Range("I136").Value = 1128.4193
$ws.Range("J136").Value = 1452.3334
$ws.Range("K136").Value = 3385.2579
$ws.Range("L136").Value = 4357.0002
$ws.Range("M136").Value = -835.2579000000001
$ws.Range("N136").Value = -9457.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 45450.78
$ws.Range("I86").Value = 72957.21000000001
$ws.Range("K86").Value = 72957.21000000001
$ws.Range("M86").Value = -71834.21000000001
$ws.Range("H89").Value = 45450.78
$ws.Range("I89").Value = 72957.21000000001
$ws.Range("K89").Value = 364786.05
$ws.Range("M89").Value = -359170.05
$ws.Range("H134").Value = 3512.3684
$ws.Range("I134").Value = 3641.5881
$ws.Range("J134").Value = 2414
$ws.Range("K134").Value = 10924.7643
$ws.Range("L134").Value = 7242
$ws.Range("M134").Value = -8389.764299999999
$ws.Range("N134").Value = -12312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("H31").Value = 22436.52
$ws.Range("I31").Value = 43203.918
$ws.Range("J31").Value = 3266.6155
$ws.Range("K31").Value = 43203.918
$ws.Range("L31").Value = 3266.6155
$ws.Range("M31").Value = -42908.918
$ws.Range("N31").Value = -3856.6155
$ws.Range("H34").Value = 22436.52
$ws.Range("I34").Value = 43203.918
$ws.Range("J34").Value = 3266.6155
$ws.Range("K34").Value = 43203.918
$ws.Range("L34").Value = 3266.6155
$ws.Range("M34").Value = -43001.918
$ws.Range("N34").Value = -3670.6155
$ws.Range("H58").Value = 25488.727
$ws.Range("I58").Value = 2262.4443
$ws.Range("J58").Value = 130007
$ws.Range("K58").Value = 2262.4443
$ws.Range("L58").Value = 130007
$ws.Range("M58").Value = -2059.4443
$ws.Range("N58").Value = -130413
$ws.Range("H132").Value = 36588764
$ws.Range("I132").Value = 35717612
$ws.Range("K132").Value = 107152836
$ws.Range("M132").Value = -107150306
$ws.Range("H134").Value = 1314
$ws.Range("I134").Value = 1359.64
$ws.Range("J134").Value = 1199.9
$ws.Range("K134").Value = 4078.92
$ws.Range("L134").Value = 3599.7
$ws.Range("M134").Value = -1543.92
$ws.Range("N134").Value = -8669.700000000001
$ws.Range("H136").Value = 25488.727
$ws.Range("I136").Value = 2262.4443
$ws.Range("J136").Value = 130007
$ws.Range("K136").Value = 6787.3329
$ws.Range("L136").Value = 390021
$ws.Range("M136").Value = -4237.3329
$ws.Range("N136").Value = -395121
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 719.5263
$ws.Range("I131").Value = 450.11765
$ws.Range("J131").Value = 778.2436
$ws.Range("K131").Value = 1350.35295
$ws.Range("L131").Value = 2334.7308
$ws.Range("M131").Value = 3689.64705
$ws.Range("N131").Value = -12414.7308

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 844635
$ws.Range("I5").Value = 2501500
$ws.Range("K5").Value = 2501500
$ws.Range("M5").Value = -2501388
$ws.Range("H132").Value = 2987.3928
$ws.Range("I132").Value = 2487.5557
$ws.Range("J132").Value = 3887.1
$ws.Range("K132").Value = 7462.6671
$ws.Range("L132").Value = 11661.3
$ws.Range("M132").Value = -4932.6671
$ws.Range("N132").Value = -16721.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 218047.81
$ws.Range("I2").Value = 444555.56
$ws.Range("J2").Value = 48167
$ws.Range("K2").Value = 444555.56
$ws.Range("L2").Value = 48167
$ws.Range("M2").Value = -444443.56
$ws.Range("N2").Value = -48391
$ws.Range("H46").Value = 1125340
$ws.Range("I46").Value = 250
$ws.Range("J46").Value = 1446794.2
$ws.Range("K46").Value = 250
$ws.Range("L46").Value = 1446794.2
$ws.Range("M46").Value = -62
$ws.Range("N46").Value = -1447170.2
$ws.Range("H55").Value = 589.28
$ws.Range("I55").Value = 315
$ws.Range("J55").Value = 718.35297
$ws.Range("K55").Value = 315
$ws.Range("L55").Value = 718.35297
$ws.Range("M55").Value = -142
$ws.Range("N55").Value = -1064.35297
$ws.Range("H132").Value = 1581.4445
$ws.Range("I132").Value = 1327.36
$ws.Range("J132").Value = 2158.9092
$ws.Range("K132").Value = 3982.08
$ws.Range("L132").Value = 6476.7276
$ws.Range("M132").Value = -1452.08
$ws.Range("N132").Value = -11536.7276
$ws.Range("H136").Value = 1575.3636
$ws.Range("I136").Value = 1485.4706
$ws.Range("J136").Value = 1881
$ws.Range("K136").Value = 4456.4118
$ws.Range("L136").Value = 5643
$ws.Range("M136").Value = -1906.4118
$ws.Range("N136").Value = -10743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 22261.111
$ws.Range("J2").Value = 38783.332
$ws.Range("L2").Value = 38783.332
$ws.Range("N2").Value = -39007.332
$ws.Range("H136").Value = 991.5
$ws.Range("I136").Value = 824.75
$ws.Range("J136").Value = 1325
$ws.Range("K136").Value = 2474.25
$ws.Range("L136").Value = 3975
$ws.Range("M136").Value = 75.75
$ws.Range("N136").Value = -9075
